$d = $word.ActiveDocument
$table = $d.Tables.Item(2)
$rowCount = $table.Rows.Count

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 1. The "180" row's constraint cell currently carries the _GoBack bookmark at the
#    end of its text. That bookmark needs to move down to the new "190" row, so
#    rewrite this paragraph without the bookmark (text/formatting unchanged).
$rfcCell = $table.Cell($rowCount - 1, 2)
$rfcP = $rfcCell.Range.Paragraphs.Item(1)
$rfcR = $rfcP.Range
$rfcXml = $xmlHeader + '<w:p w14:paraId="5F65543E" w14:textId="026B5EE4" w:rsidR="00AA75D9" w:rsidRDefault="00232D0E" w:rsidP="009C33D1"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Device data transmission shall meet RFC 1042 standard for IP datagrams</w:t></w:r></w:p>' + $xmlFooter
$rfcR.InsertXML($rfcXml)

# 2. Fill in the previously-blank "Const No." cell of the last row with "190".
$cell1 = $table.Cell($rowCount, 1)
$p1 = $cell1.Range.Paragraphs.Item(1)
$r1 = $p1.Range
$xml1 = $xmlHeader + '<w:p w14:paraId="14592956" w14:textId="77777777" w:rsidR="008B3E54" w:rsidRDefault="008B3E54" w:rsidP="009C33D1"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>190</w:t></w:r></w:p>' + $xmlFooter
$r1.InsertXML($xml1)

# 3. Fill in the previously-blank "Constraints" cell of the last row with the new
#    requirement text, carrying the _GoBack bookmark that moved out of the "180" row.
$cell2 = $table.Cell($rowCount, 2)
$p2 = $cell2.Range.Paragraphs.Item(1)
$r2 = $p2.Range
$xml2 = $xmlHeader + '<w:p w14:paraId="3B938263" w14:textId="77777777" w:rsidR="008B3E54" w:rsidRDefault="008B3E54" w:rsidP="009C33D1"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Device total cost shall not exceed $250.00</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + $xmlFooter
$r2.InsertXML($xml2)
